# Regenerate save_data: column G ("K") was recomputed using K (strikes)
# instead of the old "Strike#" metric. Write the new s_vals into column G
# for each data row (rows 2-28), leaving all other columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 5
    3  = 7
    4  = 3
    5  = 2
    6  = 6
    7  = 2
    8  = 6
    9  = 5
    10 = 5
    11 = 5
    12 = 6
    13 = 4
    14 = 5
    15 = 5
    16 = 5
    17 = 5
    18 = 4
    19 = 3
    20 = 7
    21 = 1
    22 = 6
    23 = 7
    24 = 5
    25 = 5
    26 = 4
    27 = 4
    28 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
